# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 192 (shifting the existing rows 192-241
# down to 193-242) and populate it with the new week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 192; this pushes every row that
# was at 192..241 down to 193..242, carrying along their values/formatting
# (which is exactly the data movement described by the diff).
$ws.Rows.Item(192).Insert()

# Populate the newly inserted (currently blank) row 192 with this week's data.
$ws.Cells.Item(192, 1).Value = 11
$ws.Cells.Item(192, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(192, 3).Value = "Bíobío"
$ws.Cells.Item(192, 4).Value = 44964
$ws.Cells.Item(192, 5).Value = 8
$ws.Cells.Item(192, 6).Value = "Fruta"
$ws.Cells.Item(192, 7).Value = 100108
$ws.Cells.Item(192, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(192, 9).Value = 100108005
$ws.Cells.Item(192, 10).Value = "Piña"
$ws.Cells.Item(192, 11).Value = "Caramelo"
$ws.Cells.Item(192, 12).Value = "Segunda"
$ws.Cells.Item(192, 13).Value = 150
$ws.Cells.Item(192, 14).Value = 18000
$ws.Cells.Item(192, 15).Value = 20000
$ws.Cells.Item(192, 16).Value = 19333
$ws.Cells.Item(192, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(192, 18).Value = "Ecuador"
$ws.Cells.Item(192, 19).Value = 1381
$ws.Cells.Item(192, 20).Value = 14
